# Update column F (dSF) values per repulled data / mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = -1
    3  = -3
    4  = 0
    5  = 4
    6  = -3
    7  = -4
    8  = 7
    9  = -4
    10 = -2
    12 = 1
    13 = -3
    14 = 3
    15 = -6
    16 = -2
    17 = 1
    18 = -2
    19 = -2
    20 = 4
    21 = -1
    22 = -1
    23 = -2
    24 = -1
    25 = 8
    26 = 3
    27 = -2
    28 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("F$row").Value = $newValues[$row]
}
